$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of C10: the "app-specific feedback page" link text lost its
# parenthesized URL, and "DRA sign-in page" was simplified to "sign-in page".
$ws.Range("C10").Value = "Verify that the profile fly-out should contain link to terms of use||Verify that profile fly-out should contain link to privacy statement||Verify that the profile fly-out should contain link to app-specific feedback page ||Verify that the profile fly-out should contain link to app-specific help page||Verify that the alternative profile fly-out should contain link to sign out of the platform. User returns to sign-in page."

# Row 10 got shorter (the URL was removed) so its autofit wrapped height shrank.
$ws.Rows.Item(10).RowHeight = 75

# The active selection/scroll position moved from C20/A13 to C10/A10.
$excel.Goto($ws.Range("C10"))
